# Atualizacao de bases das ligas, do dia: 17-06-2024 as 21:10
#
# The Germany Landesliga results/odds feed was re-synced: a batch of
# fixtures got their scores / odds swapped back onto the correct match
# rows (rows 12/13, 83/84, 88/90), each pair trading places in full
# (fixture id, teams, FT/HT score, result, and every odds column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value2 = 7035046
$ws.Range("E12").Value2 = "Cronenberger SC"
$ws.Range("F12").Value2 = "FC Viersen"
$ws.Range("G12").Value2 = 0
$ws.Range("H12").Value2 = 2
$ws.Range("I12").Value2 = 0
$ws.Range("L12").Value2 = 2
$ws.Range("M12").Value2 = 3.6
$ws.Range("N12").Value2 = 3
$ws.Range("O12").Value2 = 2
$ws.Range("Q12").Value2 = 3
$ws.Range("S12").Value2 = 1.8
$ws.Range("T12").Value2 = 2
$ws.Range("U12").Value2 = 2.75
$ws.Range("V12").Value2 = 1.8
$ws.Range("W12").Value2 = 2
$ws.Range("Z12").Value2 = 2
$ws.Range("AB12").Value2 = 1
$ws.Range("AC12").Value2 = -1
$ws.Range("AD12").Value2 = 1
$ws.Range("B13").Value2 = 7035047
$ws.Range("E13").Value2 = "SC Dsseldorf West"
$ws.Range("F13").Value2 = "VfL Viktoria JuchenGarzweiler"
$ws.Range("G13").Value2 = 3
$ws.Range("H13").Value2 = 4
$ws.Range("I13").Value2 = 3
$ws.Range("L13").Value2 = 1.909
$ws.Range("M13").Value2 = 3.75
$ws.Range("N13").Value2 = 3.1
$ws.Range("O13").Value2 = 2.2
$ws.Range("Q13").Value2 = 2.625
$ws.Range("S13").Value2 = 2
$ws.Range("T13").Value2 = 1.8
$ws.Range("U13").Value2 = 3
$ws.Range("V13").Value2 = 1.825
$ws.Range("W13").Value2 = 1.975
$ws.Range("Z13").Value2 = 1.625
$ws.Range("AB13").Value2 = 0.8
$ws.Range("AC13").Value2 = 0.825
$ws.Range("AD13").Value2 = -1
$ws.Range("B83").Value2 = 8075296
$ws.Range("E83").Value2 = "FC Monheim"
$ws.Range("F83").Value2 = "VFB Hilden II"
$ws.Range("G83").Value2 = 1
$ws.Range("H83").Value2 = 2
$ws.Range("J83").Value2 = 2
$ws.Range("K83").Value2 = "A"
$ws.Range("L83").Value2 = 1.533
$ws.Range("M83").Value2 = 4.75
$ws.Range("N83").Value2 = 4
$ws.Range("O83").Value2 = 1.4
$ws.Range("P83").Value2 = 5.25
$ws.Range("Q83").Value2 = 5
$ws.Range("R83").Value2 = -1.5
$ws.Range("S83").Value2 = 1.975
$ws.Range("T83").Value2 = 1.825
$ws.Range("U83").Value2 = 3.75
$ws.Range("V83").Value2 = 1.9
$ws.Range("W83").Value2 = 1.9
$ws.Range("X83").Value2 = -1
$ws.Range("Z83").Value2 = 4
$ws.Range("AA83").Value2 = -1
$ws.Range("AB83").Value2 = 0.825
$ws.Range("AD83").Value2 = 0.8999999999999999
$ws.Range("B84").Value2 = 8075530
$ws.Range("E84").Value2 = "TuRU Dsseldorf"
$ws.Range("F84").Value2 = "DV Solingen"
$ws.Range("G84").Value2 = 2
$ws.Range("H84").Value2 = 0
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = "H"
$ws.Range("L84").Value2 = 2.1
$ws.Range("M84").Value2 = 3.75
$ws.Range("N84").Value2 = 2.7
$ws.Range("O84").Value2 = 2.375
$ws.Range("P84").Value2 = 3.75
$ws.Range("Q84").Value2 = 2.45
$ws.Range("R84").Value2 = 0
$ws.Range("S84").Value2 = 1.85
$ws.Range("T84").Value2 = 1.95
$ws.Range("U84").Value2 = 3
$ws.Range("V84").Value2 = 1.85
$ws.Range("W84").Value2 = 1.95
$ws.Range("X84").Value2 = 1.375
$ws.Range("Z84").Value2 = -1
$ws.Range("AA84").Value2 = 0.8500000000000001
$ws.Range("AB84").Value2 = -1
$ws.Range("AD84").Value2 = 0.95
$ws.Range("B88").Value2 = 8076438
$ws.Range("E88").Value2 = "SV Rott 1927"
$ws.Range("F88").Value2 = "SV Breinig"
$ws.Range("G88").Value2 = 1
$ws.Range("H88").Value2 = 4
$ws.Range("I88").Value2 = 1
$ws.Range("J88").Value2 = 1
$ws.Range("K88").Value2 = "A"
$ws.Range("L88").Value2 = 1.363
$ws.Range("M88").Value2 = 5
$ws.Range("N88").Value2 = 5.5
$ws.Range("O88").Value2 = 1.65
$ws.Range("P88").Value2 = 4.5
$ws.Range("Q88").Value2 = 3.5
$ws.Range("R88").Value2 = -0.75
$ws.Range("S88").Value2 = 1.825
$ws.Range("T88").Value2 = 1.975
$ws.Range("U88").Value2 = 3.5
$ws.Range("V88").Value2 = 1.975
$ws.Range("W88").Value2 = 1.825
$ws.Range("X88").Value2 = -1
$ws.Range("Z88").Value2 = 2.5
$ws.Range("AB88").Value2 = 0.9750000000000001
$ws.Range("AC88").Value2 = 0.9750000000000001
$ws.Range("B90").Value2 = 8077795
$ws.Range("E90").Value2 = "Eiche Horn"
$ws.Range("F90").Value2 = "SVGO Bremen"
$ws.Range("G90").Value2 = 5
$ws.Range("H90").Value2 = 3
$ws.Range("I90").Value2 = 3
$ws.Range("J90").Value2 = 2
$ws.Range("K90").Value2 = "H"
$ws.Range("L90").Value2 = 1.142
$ws.Range("M90").Value2 = 7
$ws.Range("N90").Value2 = 10
$ws.Range("O90").Value2 = 1.083
$ws.Range("P90").Value2 = 11
$ws.Range("Q90").Value2 = 19
$ws.Range("R90").Value2 = -3.5
$ws.Range("S90").Value2 = 1.975
$ws.Range("T90").Value2 = 1.825
$ws.Range("U90").Value2 = 5
$ws.Range("V90").Value2 = 1.825
$ws.Range("W90").Value2 = 1.975
$ws.Range("X90").Value2 = 0.08299999999999996
$ws.Range("Z90").Value2 = -1
$ws.Range("AB90").Value2 = 0.825
$ws.Range("AC90").Value2 = 0.825
